# Camera Calibration Worksheet update
# - Archive the previous "Trial B" readings into the (previously empty) N/O/P columns
# - Promote the previous "Trial A" readings into "Trial B" (H/I/J columns)
# - Record brand new "Trial A" readings (B/C/D/E columns) from the new Waterloo trial
# - Label the new data with a note in O5
# - Bump the precision of the quadratic-fit coefficient displayed in D17:D19
# - A couple of stray manual note numbers in column R were also updated
# - Clear out a leftover "Unsure" note in F12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calibration")

# New label for the freshly-added trial
$ws.Range("O5").Value = "Waterloo trial"

# --- archive the old "Trial B" values (H8:J10) into N8:P10 (previously blank) ---
$ws.Range("N8").Value = 1.7
$ws.Range("O8").Value = 103.5
$ws.Range("P8").Value = 35

$ws.Range("N9").Value = 2.53
$ws.Range("O9").Value = 112.5
$ws.Range("P9").Value = 94.5

$ws.Range("N10").Value = 3.51
$ws.Range("O10").Value = 107.5
$ws.Range("P10").Value = 136.5

# --- promote the old "Trial A" values (B8:D12) into "Trial B" (H8:J12) ---
$ws.Range("H8").Value = 1.64
$ws.Range("I8").Value = 97
$ws.Range("J8").Value = 33.5

$ws.Range("H9").Value = 2.53
$ws.Range("I9").Value = 99.5
$ws.Range("J9").Value = 92

$ws.Range("H10").Value = 3.8
$ws.Range("I10").Value = 106
$ws.Range("J10").Value = 144

$ws.Range("H11").Value = 1.43
$ws.Range("I11").Value = 85
$ws.Range("J11").Value = 12

$ws.Range("H12").Value = 1.45
$ws.Range("I12").Value = 79.5
$ws.Range("J12").Value = 29

# --- write the brand new "Trial A" values (B8:E12) ---
$ws.Range("B8").Value = 1.6
$ws.Range("C8").Value = 99
$ws.Range("D8").Value = 23
$ws.Range("E8").Value = 1110

$ws.Range("B9").Value = 2.54
$ws.Range("C9").Value = 109
$ws.Range("D9").Value = 91
$ws.Range("E9").Value = 1340

$ws.Range("B10").Value = 3.68
$ws.Range("C10").Value = 115
$ws.Range("D10").Value = 140

$ws.Range("B11").Value = 3.48
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()

$ws.Range("B12").Value = 1.35
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("E12").Value = 1000

# Remove stray leftover note
$ws.Range("F12").ClearContents()

# Manual note numbers near the trial tables
$ws.Range("R8").Value = 5365
$ws.Range("R10").Value = -5429

# Higher-precision display format for the quadratic coefficients (C2) of the
# Trial A curve fit
$ws.Range("D17:D19").NumberFormat = "0.000000000000"
